$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.055.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.134.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.18%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.08%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.571'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.148.17'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  -3.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.383'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.688.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.95%  '

$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.193.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.143.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("E18").Value = '  -3.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '400.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.37%  '

$ws.Range("E21").Value = '  -2.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.94%  '

$ws.Range("E23").Value = '  +3.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.07%  '

$ws.Range("E26").Value = '  -0.80%  '

$ws.Range("E27").Value = '  -5.42%  '

$ws.Range("E28").Value = '  -5.25%  '

$ws.Range("E29").Value = '  -0.79%  '

$ws.Range("E30").Value = '  -1.08%  '

$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("E32").Value = '  -1.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '159.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.53%  '

$ws.Range("E35").Value = '  -1.02%  '

$ws.Range("E36").Value = '  -4.40%  '

$ws.Range("E37").Value = '  -2.16%  '

$ws.Range("E38").Value = '  -1.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.668.18'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.66'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.54%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.30'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.689'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0613'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("E46").Value = '  -2.86%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0255'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.59%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '287.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.62%  '

$ws.Range("E50").Value = '  -0.23%  '
